$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# xlPasteFormats constant
$xlPasteFormats = -4122

# Row 3: move the red "1" marker from E3 to J3
$ws.Range("E3").Copy()
$ws.Range("J3").PasteSpecial($xlPasteFormats)
$ws.Range("J3").Value = 1
$ws.Range("E3").ClearContents()

# Row 5 / Row 6: move the red "1" marker from I5 to H5, and from J6 to I6
$ws.Range("I5").Copy()
$ws.Range("H5").PasteSpecial($xlPasteFormats)
$ws.Range("H5").Value = 1

$ws.Range("G5").Copy()
$ws.Range("I5").PasteSpecial($xlPasteFormats)
$ws.Range("I5").ClearContents()

$ws.Range("J6").Copy()
$ws.Range("I6").PasteSpecial($xlPasteFormats)
$ws.Range("I6").Value = 1

$ws.Range("H6").Copy()
$ws.Range("J6").PasteSpecial($xlPasteFormats)
$ws.Range("J6").ClearContents()

$excel.CutCopyMode = $false

# Update the active selection shown in the sheet view
$ws.Range("G3").Select()
